# Apply the ImportStaff.xlsx header-row restructuring described by the diff:
#   - Sheet 1 "Staff Vitals":   A1:Z1  -> A1:AA1 (one new column; headers re-sorted,
#                                split FIRSTNAME/LASTNAME into "First Name"/"Last Name",
#                                rename "*Unique Photo ID" -> "Staff Vitals - UNIQUE_PHOTO_ID")
#   - Sheet 2 "Staff Attributes": A1:O1 -> A1:W1 (eight new columns; headers re-sorted)
#   - Sheet 3 "Staff Style":     A1:Q1 stays A1:Q1 (headers reshuffled/renamed in place)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Staff Vitals" (26 columns -> 27 columns)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Make room for the one extra header column (insert a blank column at the end
# of the current range so every existing cell keeps its style/position while
# we rewrite all the header text below).
[void]$ws1.Columns.Item(27).Insert()

$ws1.Cells.Item(1,1).Value  = "Face ID"
$ws1.Cells.Item(1,2).Value  = "First Name"
$ws1.Cells.Item(1,3).Value  = "Height"
$ws1.Cells.Item(1,4).Value  = "Last Name"
$ws1.Cells.Item(1,5).Value  = "Position"
$ws1.Cells.Item(1,6).Value  = "Salary"
$ws1.Cells.Item(1,7).Value  = "Staff Vitals - ARM_SCALE"
$ws1.Cells.Item(1,8).Value  = "Staff Vitals - BODYLENGTH"
$ws1.Cells.Item(1,9).Value  = "Staff Vitals - BODY_SHAPE"
$ws1.Cells.Item(1,10).Value = "Staff Vitals - CURRENT_TEAM"
$ws1.Cells.Item(1,11).Value = "Staff Vitals - EYE_COLOR"
$ws1.Cells.Item(1,12).Value = "Staff Vitals - GENDER"
$ws1.Cells.Item(1,13).Value = "Staff Vitals - HAIR_LENGTH"
$ws1.Cells.Item(1,14).Value = "Staff Vitals - HAND_SCALE"
$ws1.Cells.Item(1,15).Value = "Staff Vitals - HEIGHT_CM"
$ws1.Cells.Item(1,16).Value = "Staff Vitals - LOWER_SCALE"
$ws1.Cells.Item(1,17).Value = "Staff Vitals - NECK_HEAD_SCALE"
$ws1.Cells.Item(1,18).Value = "Staff Vitals - PERSONALITY"
$ws1.Cells.Item(1,19).Value = "Staff Vitals - POSITION"
$ws1.Cells.Item(1,20).Value = "Staff Vitals - SALARY"
$ws1.Cells.Item(1,21).Value = "Staff Vitals - SHOULDERWIDTH"
$ws1.Cells.Item(1,22).Value = "Staff Vitals - SKINCOLOR"
$ws1.Cells.Item(1,23).Value = "Staff Vitals - SKINTYPE"
$ws1.Cells.Item(1,24).Value = "Staff Vitals - UNIQUE_PHOTO_ID"
$ws1.Cells.Item(1,25).Value = "Staff Vitals - WINGSPAN_CM"
$ws1.Cells.Item(1,26).Value = "Staff Vitals - YEARS_IN_LEAGUE"
$ws1.Cells.Item(1,27).Value = "Staff Vitals - YEARS_LEFT"

# ---------------------------------------------------------------------------
# Sheet 2: "Staff Attributes" (15 columns -> 23 columns)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Insert eight blank columns right after column A so each new cell inherits
# the header style (bold/centered/bordered) from its left neighbour, then
# shift the rest of the existing header cells out of the way.
for ($i = 0; $i -lt 8; $i++) {
    [void]$ws2.Columns.Item(2).Insert()
}

$ws2.Cells.Item(1,1).Value  = "Analytics"
$ws2.Cells.Item(1,2).Value  = "Business"
$ws2.Cells.Item(1,3).Value  = "Charisma"
$ws2.Cells.Item(1,4).Value  = "Current Team"
$ws2.Cells.Item(1,5).Value  = "Defensive Coaching"
$ws2.Cells.Item(1,6).Value  = "Offensive Coaching"
$ws2.Cells.Item(1,7).Value  = "Potential"
$ws2.Cells.Item(1,8).Value  = "Sports Medicine"
$ws2.Cells.Item(1,9).Value  = "Staff Attributes - BUSINESS"
$ws2.Cells.Item(1,10).Value = "Staff Attributes - CONTRACTS"
$ws2.Cells.Item(1,11).Value = "Staff Attributes - DEFENSE"
$ws2.Cells.Item(1,12).Value = "Staff Attributes - MAX_BUSINESS"
$ws2.Cells.Item(1,13).Value = "Staff Attributes - MAX_CONTRACTS"
$ws2.Cells.Item(1,14).Value = "Staff Attributes - MAX_DEFENSE"
$ws2.Cells.Item(1,15).Value = "Staff Attributes - MAX_OFFENSE"
$ws2.Cells.Item(1,16).Value = "Staff Attributes - MAX_SCOUTING"
$ws2.Cells.Item(1,17).Value = "Staff Attributes - MAX_TRADING"
$ws2.Cells.Item(1,18).Value = "Staff Attributes - MAX_TRAINING"
$ws2.Cells.Item(1,19).Value = "Staff Attributes - OFFENSE"
$ws2.Cells.Item(1,20).Value = "Staff Attributes - POTENTIAL"
$ws2.Cells.Item(1,21).Value = "Staff Attributes - SCOUTING"
$ws2.Cells.Item(1,22).Value = "Staff Attributes - TRADING"
$ws2.Cells.Item(1,23).Value = "Staff Attributes - TRAINING"

# ---------------------------------------------------------------------------
# Sheet 3: "Staff Style" (17 columns, same count - headers reshuffled/renamed)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Cells.Item(1,1).Value  = "Balanced Proficiency"
$ws3.Cells.Item(1,2).Value  = "Defense Proficiency"
$ws3.Cells.Item(1,3).Value  = "Grit & Grind Proficiency"
$ws3.Cells.Item(1,4).Value  = "Pace &Space Proficiency"
$ws3.Cells.Item(1,5).Value  = "Perimeter Centric Proficiency"
$ws3.Cells.Item(1,6).Value  = "Seven Seconds Proficiency"
$ws3.Cells.Item(1,7).Value  = "Staff Style - ACTIVE_SYSTEM"
$ws3.Cells.Item(1,8).Value  = "Staff Style - GUARDS_VS_FORWARDS"
$ws3.Cells.Item(1,9).Value  = "Staff Style - INSIDE_VS_OUTSIDE"
$ws3.Cells.Item(1,10).Value = "Staff Style - OFFENSE_VS_DEFENSE"
$ws3.Cells.Item(1,11).Value = "Staff Style - PERIMETER_CENTERIC_PROFICIENCY"
$ws3.Cells.Item(1,12).Value = "Staff Style - POST_CENTRIC_PROFICIENCY"
$ws3.Cells.Item(1,13).Value = "Staff Style - PREFERRED_SYSTEM"
$ws3.Cells.Item(1,14).Value = "Staff Style - STYLE_N#1"
$ws3.Cells.Item(1,15).Value = "Staff Style - STYLE_N#2"
$ws3.Cells.Item(1,16).Value = "Staff Style - STYLE_N#3"
$ws3.Cells.Item(1,17).Value = "Triangle Proficiency"

Write-Output "ImportStaff header rows updated."
